$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the second data row (row 3) with a new policy number (NroPoliza, col E)
# and a new claim date (FechaSiniestro, col G) so validation can be re-tested
# against a fresh record that still needs phone/email loaded.
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "11111003199 "

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "31/05/2021"

# Move the active selection to G4, matching where the user left off.
$ws.Range("G4").Select()
